# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2"  = 408
    "F4"  = 6801
    "F5"  = 373
    "F7"  = 3419
    "F8"  = 35
    "F9"  = 25
    "F10" = 48
    "F11" = 848
    "F13" = 5411
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
